# Rename sheet3: 部族別_代碼 -> 部落+族別_代碼
$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "部落+族別_代碼"

# New sheet: 族語_代碼 (inserted right after 部落+族別_代碼)
$ws4 = $wb.Worksheets.Add($null, $ws3)
$ws4.Name = "族語_代碼"

$s4data = @(
    @("族語", "族語_代碼"),
    @("阿美族語", "'01"),
    @("泰雅族語", "'02"),
    @("排灣族語", "'03"),
    @("布農族語", "'04"),
    @("卑南族語", "'05"),
    @("魯凱族語", "'06"),
    @("鄒族語", "'07"),
    @("賽夏族語", "'08"),
    @("雅美族語", "'09"),
    @("邵族語", "'10"),
    @("噶瑪蘭族語", "'11"),
    @("太魯閣族語", "'12"),
    @("撒奇萊雅族語", "'13"),
    @("賽德克族語", "'14"),
    @("拉阿魯哇族語", "'15"),
    @("卡那卡那富族語", "'16"),
)
for ($i = 0; $i -lt $s4data.Count; $i++) {
    $r = $i + 1
    $row = $s4data[$i]
    $ws4.Cells.Item($r, 1).Value = $row[0]
    $ws4.Cells.Item($r, 2).Value = $row[1]
}
$ws4.Range("A1:B1").Font.Bold = $true
$ws4.Range("A1:B1").HorizontalAlignment = -4108
$ws4.Range("A1:B1").VerticalAlignment = -4160
$ws4.Range("A1:B1").Borders.LineStyle = 1

# New sheet: 族語+方言_代碼 (inserted right after 族語_代碼)
$ws5 = $wb.Worksheets.Add($null, $ws4)
$ws5.Name = "族語+方言_代碼"

$ws5.Cells.Item(1, 1).Value = "族語"
$ws5.Cells.Item(1, 2).Value = "方言`n原住民族語言能力認證測驗"
$ws5.Cells.Item(1, 3).Value = "方言別_代碼"
$ws5.Cells.Item(1, 4).Value = "族語+方言_代碼"
$ws5.Range("A1:D1").Font.Bold = $true
$ws5.Range("A1:D1").HorizontalAlignment = -4108
$ws5.Range("A1:D1").VerticalAlignment = -4160
$ws5.Range("A1:D1").Borders.LineStyle = 1

$s5data = @(
    @("阿美族語", "南勢阿美語", "'01", "'0101"),
    @("阿美族語", "秀姑巒阿美語", "'02", "'0102"),
    @("阿美族語", "海岸阿美語", "'03", "'0103"),
    @("阿美族語", "馬蘭阿美語", "'04", "'0104"),
    @("阿美族語", "恆春阿美語", "'05", "'0105"),
    @("泰雅族語", "賽考利克泰雅語", "'01", "'0201"),
    @("泰雅族語", "四季泰雅語", "'02", "'0202"),
    @("泰雅族語", "澤敖利泰雅語", "'03", "'0203"),
    @("泰雅族語", "宜蘭澤敖利泰雅語", "'04", "'0204"),
    @("泰雅族語", "汶水泰雅語", "'05", "'0205"),
    @("泰雅族語", "萬大泰雅語", "'06", "'0206"),
    @("排灣族語", "東排灣語", "'01", "'0301"),
    @("排灣族語", "北排灣語", "'02", "'0302"),
    @("排灣族語", "中排灣語", "'03", "'0303"),
    @("排灣族語", "南排灣語", "'04", "'0304"),
    @("布農族語", "卓群布農語", "'01", "'0401"),
    @("布農族語", "卡群布農語", "'02", "'0402"),
    @("布農族語", "丹群布農語", "'03", "'0403"),
    @("布農族語", "巒群布農語", "'04", "'0404"),
    @("布農族語", "郡群布農語", "'05", "'0405"),
    @("卑南族語", "知本卑南語", "'01", "'0501"),
    @("卑南族語", "南王卑南語", "'02", "'0502"),
    @("卑南族語", "西群卑南語", "'03", "'0503"),
    @("卑南族語", "建和卑南語", "'04", "'0504"),
    @("魯凱族語", "霧臺魯凱語", "'01", "'0601"),
    @("魯凱族語", "東魯凱語", "'02", "'0602"),
    @("魯凱族語", "多納魯凱語", "'03", "'0603"),
    @("魯凱族語", "萬山魯凱語", "'04", "'0604"),
    @("魯凱族語", "茂林魯凱語", "'05", "'0605"),
    @("魯凱族語", "大武魯凱語", "'06", "'0606"),
    @("鄒族語", "鄒語", "'01", "'0701"),
    @("賽夏族語", "賽夏語", "'01", "'0801"),
    @("雅美族語", "雅美語", "'01", "'0901"),
    @("邵族語", "邵語", "'01", "'1001"),
    @("噶瑪蘭族語", "噶瑪蘭語", "'01", "'1101"),
    @("太魯閣族語", "太魯閣語", "'01", "'1201"),
    @("撒奇萊雅族語", "撒奇萊雅語", "'01", "'1301"),
    @("賽德克族語", "都達語", "'01", "'1401"),
    @("賽德克族語", "德固達雅語", "'02", "'1402"),
    @("賽德克族語", "德路固語", "'03", "'1403"),
    @("拉阿魯哇族語", "拉阿魯哇語", "'01", "'1501"),
    @("卡那卡那富族語", "卡那卡那富語", "'01", "'1601"),
)
for ($i = 0; $i -lt $s5data.Count; $i++) {
    $r = $i + 2
    $row = $s5data[$i]
    $ws5.Cells.Item($r, 1).Value = $row[0]
    $ws5.Cells.Item($r, 2).Value = $row[1]
    $ws5.Cells.Item($r, 3).Value = $row[2]
    $ws5.Cells.Item($r, 4).Value = $row[3]
}
